# Start scene to test Arthur out
#
# - Adds four new task cells for the upcoming "Arthur" work:
#     I7 (Tuesday column)   "Add song to scene"
#     J2 (Wednesday column) "Pseudocode out Arthur Jump"
#     J3 (Wednesday column) "Implement Arthur Jump"
#     J4 (Wednesday column) "Implement Arthur Sprite"
# - Gives I5 / I6 the same yellow-fill look as the rest of the Tuesday
#   column (I2/I3).
# - Widens column J very slightly.
# - Leaves the active selection on I5.
#
# Cell formats are copied from an existing cell that already has the
# desired look (PasteSpecial -> formats only) instead of being built up
# property-by-property; doing it property-by-property can leave a
# transient, unused cell style registered in the workbook because each
# property write commits its own style lookup/allocation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cells ---------------------------------------------------------
# Values are written in the same order the corresponding shared strings
# appear in the target workbook (I7, then J2, J3, J4).

$ws.Range("I4").Copy()
$ws.Range("I7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I7").Value = "Add song to scene"

$ws.Range("I5").Copy()
$ws.Range("J2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("J2").Value = "Pseudocode out Arthur Jump"

$ws.Range("J3").Value = "Implement Arthur Jump"
$ws.Range("J4").Value = "Implement Arthur Sprite"

# --- Restyle existing cells ---------------------------------------------
# I5 / I6 pick up the yellow fill already used across the Tuesday column.

$ws.Range("I2").Copy()
$ws.Range("I5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I6").PasteSpecial(-4122)   # xlPasteFormats

# --- Column width --------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 21.3

# --- Selection -------------------------------------------------------
[void]$ws.Range("I5").Select()
